$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Paragraph 3 ("Micka Kabongo ..."): replace " -----" with " " + tab + "260800865" ---
$p3 = $d.Range(41, 47)
$xml3 = '<w:p ' + $wNs + '><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:tab/></w:r><w:r><w:t>260800865</w:t></w:r></w:p>'
$p3.InsertXML($xml3)

# --- Paragraph 2 ("Esa Khan 260706611"): split into "Esa Khan " + 3 tabs + bookmark + "260706611" ---
$p2 = $d.Range(9, 27)
$xml2 = '<w:p ' + $wNs + '><w:r><w:t xml:space="preserve">Esa Khan </w:t></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:r><w:tab/></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>260706611</w:t></w:r></w:p>'
$p2.InsertXML($xml2)

# --- Paragraph 1 ("Group 29"): remove the now-duplicate _GoBack bookmark ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

Write-Output $d.Content.Text
